$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '61.899.02'
$ws.Range('E2').Value = '  -0.97%  '
Set-TextValue 'D3' '2.897.98'
$ws.Range('E3').Value = '  -1.68%  '
$ws.Range('E4').Value = '  +0.08%  '
Set-TextValue 'D5' '568.73'
$ws.Range('E5').Value = '  -3.32%  '
Set-TextValue 'D6' '143.98'
$ws.Range('E6').Value = '  -1.49%  '
$ws.Range('E7').Value = '  -0.03%  '
Set-TextValue 'D8' '0.502'
$ws.Range('E8').Value = '  -0.95%  '
Set-TextValue 'D9' '2.896.10'
$ws.Range('E9').Value = '  -1.70%  '
Set-TextValue 'D10' '6.90'
$ws.Range('E10').Value = '  -1.89%  '
$ws.Range('E11').Value = '  -2.28%  '
$ws.Range('E12').Value = '  -1.17%  '
Set-TextValue 'D13' '0.0000230'
$ws.Range('E13').Value = '  -1.46%  '
Set-TextValue 'D14' '32.38'
$ws.Range('E14').Value = '  +0.38%  '
Set-TextValue 'D15' '0.125'
Set-TextValue 'D16' '3.378.05'
$ws.Range('E16').Value = '  -1.68%  '
Set-TextValue 'D17' '61.865.31'
$ws.Range('E17').Value = '  -0.95%  '
$ws.Range('E18').Value = '  -1.74%  '
Set-TextValue 'D19' '2.893.28'
$ws.Range('E19').Value = '  -1.74%  '
Set-TextValue 'D20' '431.33'
$ws.Range('E20').Value = '  -0.63%  '
$ws.Range('E21').Value = '  -3.49%  '
Set-TextValue 'D22' '0.652'
$ws.Range('E22').Value = '  -1.44%  '
$ws.Range('E23').Value = '  -1.40%  '
Set-TextValue 'D24' '78.84'
$ws.Range('E24').Value = '  -1.54%  '
Set-TextValue 'D25' '12.07'
$ws.Range('E25').Value = '  +1.06%  '
Set-TextValue 'D26' '10.07'
$ws.Range('E26').Value = '  -8.96%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('E29').Value = '  +9.89%  '
Set-TextValue 'D30' '7.01'
$ws.Range('E30').Value = '  -2.25%  '
$ws.Range('E31').Value = '  -2.86%  '
$ws.Range('E32').Value = '  -5.69%  '
Set-TextValue 'D33' '1.00'
$ws.Range('E33').Value = '  +0.20%  '
$ws.Range('E34').Value = '  -1.46%  '
Set-TextValue 'D35' '25.56'
$ws.Range('E35').Value = '  -2.37%  '
$ws.Range('E36').Value = '  -3.49%  '
$ws.Range('E37').Value = '  -3.38%  '
Set-TextValue 'D38' '48.82'
Set-TextValue 'D39' '2.85'
$ws.Range('E39').Value = '  -5.23%  '
$ws.Range('E40').Value = '  -4.29%  '
$ws.Range('E41').Value = '  +0.23%  '
Set-TextValue 'D42' '8.14'
Set-TextValue 'D43' '40.14'
$ws.Range('E43').Value = '  +3.89%  '
$ws.Range('E44').Value = '  -1.99%  '
Set-TextValue 'D45' '2.697.76'
$ws.Range('E45').Value = '  +0.28%  '
$ws.Range('E46').Value = '  -0.49%  '
Set-TextValue 'D47' '131.80'
$ws.Range('E47').Value = '  -2.45%  '
Set-TextValue 'D48' '347.19'
$ws.Range('E48').Value = '  -2.12%  '
$ws.Range('E50').Value = '  -1.33%  '
